# Update "想去人数" (interested-count) figures in both the "展览" sheet
# (sheet 1) and the "全部类型" sheet (sheet 4) to reflect the refreshed
# scrape output.

$wb = $excel.ActiveWorkbook

# Sheet 1 ("展览") - column F holds the counts for rows 5, 8, 9, 11, 12, 13
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F5").Value = 5417
$wsExpo.Range("F8").Value = 135
$wsExpo.Range("F9").Value = 2409
$wsExpo.Range("F11").Value = 56
$wsExpo.Range("F12").Value = 2260
$wsExpo.Range("F13").Value = 83

# Sheet 4 ("全部类型") - same events, rows 5, 10, 11, 14, 15, 16
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F5").Value = 5417
$wsAll.Range("F10").Value = 135
$wsAll.Range("F11").Value = 2409
$wsAll.Range("F14").Value = 56
$wsAll.Range("F15").Value = 2260
$wsAll.Range("F16").Value = 83
